# Logic tree input file updated
#
# Inserts two new rows into the decision-tree table on Sheet1:
#   - a new "Possible_Problem" leaf row right after the
#     "driving conditions" question row (old row 6 -> becomes new row 7)
#   - a new "Possible_Problem" leaf row right after the
#     "white smoke from tailpipe" question row (old row 9 -> becomes new row 14)
# Both new rows reuse the exact same Node2/Relationship/Possible_Problem
# text that already appears for the "coolant level" branch (row 4), just
# with a different Node1 (question) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$possibleProblem = "Possible_Problem"
$possibleProblemText = "Possible_Problem:25% Low Coolant level`n20% Engine Cooling Fan Failure`n15% Restricted radiator airflow`n15% Restricted radiator coolant flow`n10% Thermostat Failure`n10% Cylinder head gasket`n5% Waterpump Impeller"

$drivingConditionsQuestion = "Problem:What are the driving conditions when the overheating occurs,Freeway Driving Only, Stop & Go traffic Only, All of the time? (Please answer as: Freeway driving, Stop & Go, All of the time)"
$whiteSmokeQuestion = "Problem:With the engine running, is there a large amount of white smoke coming from the tailpipe? (Please answer as: Yes, No, Not Sure)"

# --- Insert new row before old row 7 (pushes old rows 7..15 to 8..16) ---
$ws.Rows("7:7").Insert()
$ws.Rows("7:7").RowHeight = 409.6
$ws.Range("A7").Value = $drivingConditionsQuestion
$ws.Range("B7").Value = $possibleProblem
$ws.Range("C7").Value = $possibleProblemText
$ws.Range("C7").WrapText = $true

# --- Insert new row before old row 13, now at row 14 after the first shift
#     (pushes current rows 14..16 to 15..17) ---
$ws.Rows("14:14").Insert()
$ws.Rows("14:14").RowHeight = 409.6
$ws.Range("A14").Value = $whiteSmokeQuestion
$ws.Range("B14").Value = $possibleProblem
$ws.Range("C14").Value = $possibleProblemText
$ws.Range("C14").WrapText = $true

# --- Restore the view state: scrolled down near the bottom of the table
#     with C14 (the newly added cell) selected ---
$ws.Range("C14").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "Inserted 2 rows; dimension now A1:C17"
